$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF") need the same bold/centered/
# bordered header style already used by the other header cells (e.g. H1).
# Copy H1's formatting into I1:J1 first, then set the header text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill in the new data columns I (I0) and J (IF) for rows 2-12.
$iValues = @(8, 9, 8, 8, 5, 8, 8, 8, 7, 7, 8)
$jValues = @(9, 9, 9, 8, 6, 9, 9, 8, 8, 7, 8)

for ($r = 2; $r -le 12; $r++) {
    $ws.Cells.Item($r, 9).Value = $iValues[$r - 2]
    $ws.Cells.Item($r, 10).Value = $jValues[$r - 2]
}
